$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I and J, matching style of existing headers (H1 etc.)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I values per row (multiplier), J = H * I (computed result)
$iValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 6
    11 = 7
    12 = 7
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
}

for ($r = 2; $r -le 19; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    $i = $iValues[$r]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $h * $i
}
